# Edit script: restore/normalize ZBRA price-history rows that had been
# accidentally overwritten with data from other tickers (TRI, G, CRM, ...).
# For every data row except the one that already holds the correct ZBRA
# values (row 23), we:
#   - set open/close/high/low (D/E/F/G) to the correct ZBRA price values
#   - set shares_outstanding (H) to the correct ZBRA share count (50845151)
#   - set fixed_ticker (I) back to "ZBRA"
# Once every row points at the "ZBRA" shared string, the now-unused ticker
# strings (TRI, G, CRM, NOW, META, ...) are dropped from sharedStrings.xml
# automatically on save, shrinking uniqueCount from 81 to 44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = 90.62999725341795; E = 92.08000183105467; F = 98.29000091552734; G = 87.98000335693359; H = 50845151 },
    @{ Row = 3; D = 112.620002746582; E = 107.629997253418; F = 117; G = 103.5; H = 50845151 },
    @{ Row = 4; D = 76.19999694824219; E = 76.90000152587891; F = 80.83000183105469; G = 71.20999908447266; H = 50845151 },
    @{ Row = 5; D = 67.68000030517578; E = 60.40000152587891; F = 67.68000030517578; G = 54.88999938964844; H = 50845151 },
    @{ Row = 6; D = 68.01000213623047; E = 62.56000137329102; F = 68.48999786376953; G = 62.04999923706055; H = 50845151 },
    @{ Row = 7; D = 50.0099983215332; E = 53.0099983215332; F = 53.72000122070312; G = 46.13000106811523; H = 50845151 },
    @{ Row = 8; D = 69.13999938964844; E = 65.83999633789062; F = 70.94999694824219; G = 63.84999847412109; H = 50845151 },
    @{ Row = 9; D = 86.58000183105469; E = 83.66999816894531; F = 88.02999877929688; G = 81.62999725341797; H = 50845151 },
    @{ Row = 10; D = 91.16000366210938; E = 94.26999664306641; F = 96.54000091552734; G = 86.81999969482422; H = 50845151 },
    @{ Row = 11; D = 100.9100036621094; E = 101.7200012207031; F = 103.6900024414062; G = 97.83000183105467; H = 50845151 },
    @{ Row = 12; D = 108.4800033569336; E = 115.9899978637695; F = 117.4400024414062; G = 108.4100036621094; H = 50845151 },
    @{ Row = 13; D = 104.0999984741211; E = 123.1600036621094; F = 126.4899978637695; G = 102.75; H = 50845151 },
    @{ Row = 14; D = 137.8600006103516; E = 134.8300018310547; F = 146.6999969482422; G = 134.2700042724609; H = 50845151 },
    @{ Row = 15; D = 141.8899993896484; E = 137.9299926757812; F = 150.6600036621094; G = 136.1600036621094; H = 50845151 },
    @{ Row = 16; D = 177.7200012207031; E = 166.3000030517578; F = 177.8000030517578; G = 147.8899993896484; H = 50845151 },
    @{ Row = 17; D = 155.9900054931641; E = 173.6000061035156; F = 179.8099975585938; G = 146.5200042724609; H = 50845151 },
    @{ Row = 18; D = 211.2400054931641; E = 211.1399993896484; F = 237.1499938964844; G = 209.5; H = 50845151 },
    @{ Row = 19; D = 214.0899963378907; E = 210.8899993896484; F = 220.4600067138672; G = 177.0599975585938; H = 50845151 },
    @{ Row = 20; D = 207.6100006103516; E = 237.8699951171875; F = 241.759994506836; G = 188.8000030517578; H = 50845151 },
    @{ Row = 21; D = 259.0299987792969; E = 239.0200042724609; F = 260.3999938964844; G = 238.8600006103516; H = 50845151 },
    @{ Row = 22; D = 177.5599975585938; E = 229.6600036621093; F = 248.2100067138672; G = 170.6300048828125; H = 50845151 },
    @{ Row = 24; D = 256.8099975585937; E = 283.6400146484375; F = 304.5700073242188; G = 252.5200042724609; H = 50845151 },
    @{ Row = 25; D = 385.75; E = 387.8299865722656; F = 419.2099914550781; G = 372.4599914550781; H = 50845151 },
    @{ Row = 26; D = 489.9200134277344; E = 487.739990234375; F = 518.6599731445312; G = 483.2900085449219; H = 50845151 },
    @{ Row = 27; D = 530; E = 552.47998046875; F = 554.219970703125; G = 502.6300048828125; H = 50845151 },
    @{ Row = 28; D = 517.6199951171875; E = 533.9500122070312; F = 543.489990234375; G = 491.989990234375; H = 50845151 },
    @{ Row = 29; D = 592.0800170898438; E = 509.1199951171875; F = 599.72998046875; G = 462.3200073242188; H = 50845151 },
    @{ Row = 30; D = 414.6300048828125; E = 369.6600036621094; F = 437.2999877929688; G = 368.3500061035156; H = 50845151 },
    @{ Row = 31; D = 292.2300109863281; E = 357.6900024414062; F = 358.3900146484375; G = 283.7200012207031; H = 50845151 },
    @{ Row = 32; D = 264.9400024414062; E = 283.2200012207031; F = 288.6099853515625; G = 248.1699981689453; H = 50845151 },
    @{ Row = 33; D = 260.8699951171875; E = 316.1799926757812; F = 320.1199951171875; G = 254.6900024414062; H = 50845151 },
    @{ Row = 34; D = 315.1300048828125; E = 288.0299987792969; F = 316.2300109863281; G = 278.2099914550781; H = 50845151 },
    @{ Row = 35; D = 293.8900146484375; E = 307.9599914550781; F = 320.5499877929688; G = 283.5400085449219; H = 50845151 },
    @{ Row = 36; D = 236.5299987792969; E = 209.4299926757812; F = 236.6999969482422; G = 196.1300048828125; H = 50845151 },
    @{ Row = 37; D = 268.5599975585937; E = 239.5500030517578; F = 270.739990234375; G = 236.1900024414062; H = 50845151 },
    @{ Row = 38; D = 302.5899963378906; E = 314.5599975585937; F = 322.9500122070312; G = 266.75; H = 50845151 },
    @{ Row = 39; D = 308.8200073242188; E = 351.1900024414062; F = 372.9299926757813; G = 303.3599853515625; H = 50845151 },
    @{ Row = 40; D = 369.8399963378906; E = 381.9700012207031; F = 394.5700073242188; G = 359.3699951171875; H = 50845151 },
    @{ Row = 41; D = 386.5; E = 391.9400024414063; F = 427.760009765625; G = 376.4299926757813; H = 50845151 },
    @{ Row = 42; D = 280.9599914550781; E = 250.3200073242188; F = 289.9299926757812; G = 205.729995727539; H = 50845151 },
    @{ Row = 43; D = 307.5400085449219; E = 339.0199890136719; F = 352.6600036621094; G = 307.5400085449219; H = 50845151 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = "ZBRA"
}
